# Auto-generated edit script applying numeric updates per the diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 245.8
$ws.Range("I2").Value = 257.5
$ws.Range("J2").Value = 199
$ws.Range("K2").Value = 257.5
$ws.Range("L2").Value = 199
$ws.Range("M2").Value = -144.5
$ws.Range("N2").Value = -425
$ws.Range("H28").Value = 27027750
$ws.Range("I28").Value = 34483404
$ws.Range("J28").Value = 1002.125
$ws.Range("K28").Value = 34483404
$ws.Range("L28").Value = 1002.125
$ws.Range("M28").Value = -34482919
$ws.Range("N28").Value = -1972.125
$ws.Range("I62").Value = 1061.75
$ws.Range("J62").Value = 937.5
$ws.Range("K62").Value = 1061.75
$ws.Range("L62").Value = 937.5
$ws.Range("M62").Value = -437.75
$ws.Range("N62").Value = -2185.5
$ws.Range("I65").Value = 1061.75
$ws.Range("J65").Value = 937.5
$ws.Range("K65").Value = 5308.75
$ws.Range("L65").Value = 4687.5
$ws.Range("M65").Value = -2188.75
$ws.Range("N65").Value = -10927.5
$ws.Range("H86").Value = 2825.0667
$ws.Range("I86").Value = 3243.6667
$ws.Range("K86").Value = 3243.6667
$ws.Range("M86").Value = -2120.6667
$ws.Range("H89").Value = 2825.0667
$ws.Range("I89").Value = 3243.6667
$ws.Range("K89").Value = 16218.3335
$ws.Range("M89").Value = -10602.3335
$ws.Range("H100").Value = 4443.8945
$ws.Range("I100").Value = 2316.7144
$ws.Range("J100").Value = 10400
$ws.Range("K100").Value = 2316.7144
$ws.Range("L100").Value = 10400
$ws.Range("M100").Value = -1775.7144
$ws.Range("N100").Value = -11482
$ws.Range("H135").Value = 14287164
$ws.Range("I135").Value = 14707314
$ws.Range("J135").Value = 2036
$ws.Range("K135").Value = 132365826
$ws.Range("L135").Value = 18324
$ws.Range("M135").Value = -132363291
$ws.Range("N135").Value = -23394
$ws.Range("H137").Value = 5131.875
$ws.Range("I137").Value = 4722.143
$ws.Range("K137").Value = 14166.429
$ws.Range("M137").Value = -11616.429
$ws.Range("H138").Value = 2698.9714
$ws.Range("I138").Value = 1302.85
$ws.Range("J138").Value = 3257.42
$ws.Range("K138").Value = 3908.55
$ws.Range("L138").Value = 9772.26
$ws.Range("M138").Value = 1231.45
$ws.Range("N138").Value = -20052.26
$ws.Range("H141").Value = 2323.3865
$ws.Range("I141").Value = 1877.8049
$ws.Range("K141").Value = 5633.4147
$ws.Range("M141").Value = -453.4147000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H32").Value = 10757676
$ws.Range("J32").Value = 11599.8
$ws.Range("L32").Value = 11599.8
$ws.Range("N32").Value = -12173.8
$ws.Range("H61").Value = 2483.0894
$ws.Range("I61").Value = 1955.2916
$ws.Range("K61").Value = 1955.2916
$ws.Range("M61").Value = -1743.2916
$ws.Range("H74").Value = 1394.7273
$ws.Range("I74").Value = 1394.7273
$ws.Range("K74").Value = 1394.7273
$ws.Range("M74").Value = -520.7273
$ws.Range("H77").Value = 1394.7273
$ws.Range("I77").Value = 1394.7273
$ws.Range("K77").Value = 6973.636500000001
$ws.Range("M77").Value = -2605.636500000001
$ws.Range("H110").Value = 896.95654
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 2597.275
$ws.Range("I132").Value = 1749.3572
$ws.Range("K132").Value = 5248.071599999999
$ws.Range("M132").Value = -2718.071599999999
$ws.Range("H136").Value = 2483.0894
$ws.Range("I136").Value = 1955.2916
$ws.Range("K136").Value = 5865.8748
$ws.Range("M136").Value = -3315.8748

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H22").Value = 2159.3
$ws.Range("J22").Value = 68.333336
$ws.Range("L22").Value = 68.333336
$ws.Range("N22").Value = -414.333336
$ws.Range("H86").Value = 2881.5
$ws.Range("I86").Value = 1332.8889
$ws.Range("K86").Value = 1332.8889
$ws.Range("M86").Value = -209.8888999999999
$ws.Range("H89").Value = 2881.5
$ws.Range("I89").Value = 1332.8889
$ws.Range("K89").Value = 6664.4445
$ws.Range("M89").Value = -1048.4445
$ws.Range("H107").Value = 800
$ws.Range("J107").Value = 800
$ws.Range("L107").Value = 800
$ws.Range("N107").Value = -4640
$ws.Range("H134").Value = 4650.2856
$ws.Range("I134").Value = 4615.846
$ws.Range("J134").Value = 4670.636
$ws.Range("K134").Value = 13847.538
$ws.Range("L134").Value = 14011.908
$ws.Range("M134").Value = -11312.538
$ws.Range("N134").Value = -19081.908

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5734.625
$ws.Range("I31").Value = 4375.6
$ws.Range("J31").Value = 7999.6665
$ws.Range("K31").Value = 4375.6
$ws.Range("L31").Value = 7999.6665
$ws.Range("M31").Value = -4080.6
$ws.Range("N31").Value = -8589.666499999999
$ws.Range("H34").Value = 5734.625
$ws.Range("I34").Value = 4375.6
$ws.Range("J34").Value = 7999.6665
$ws.Range("K34").Value = 4375.6
$ws.Range("L34").Value = 7999.6665
$ws.Range("M34").Value = -4173.6
$ws.Range("N34").Value = -8403.666499999999
$ws.Range("H132").Value = 1755.8889
$ws.Range("I132").Value = 1340.8541
$ws.Range("K132").Value = 4022.5623
$ws.Range("M132").Value = -1492.5623
$ws.Range("H134").Value = 3149.6438
$ws.Range("I134").Value = 2038.1923
$ws.Range("J134").Value = 5901.8096
$ws.Range("K134").Value = 6114.5769
$ws.Range("L134").Value = 17705.4288
$ws.Range("M134").Value = -3579.5769
$ws.Range("N134").Value = -22775.4288

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1194.3
$ws.Range("I68").Value = 1451.6
$ws.Range("J68").Value = 937
$ws.Range("K68").Value = 4354.799999999999
$ws.Range("L68").Value = 2811
$ws.Range("M68").Value = -3543.799999999999
$ws.Range("N68").Value = -4433
$ws.Range("H71").Value = 1194.3
$ws.Range("I71").Value = 1451.6
$ws.Range("J71").Value = 937
$ws.Range("K71").Value = 13064.4
$ws.Range("L71").Value = 8433
$ws.Range("M71").Value = -9008.4
$ws.Range("N71").Value = -16545
$ws.Range("H132").Value = 1269.4
$ws.Range("J132").Value = 1199.25
$ws.Range("L132").Value = 10793.25
$ws.Range("N132").Value = -15853.25
$ws.Range("H133").Value = 4945
$ws.Range("I133").Value = 4945
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 14835
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -9775
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 6392.5454
$ws.Range("I134").Value = 1414.75
$ws.Range("K134").Value = 4244.25
$ws.Range("M134").Value = 825.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2773.1428
$ws.Range("I102").Value = 1690.7059
$ws.Range("K102").Value = 1690.7059
$ws.Range("M102").Value = -68.70589999999993
$ws.Range("H132").Value = 3275.0386
$ws.Range("I132").Value = 3127.0833
$ws.Range("J132").Value = 3401.8572
$ws.Range("K132").Value = 9381.249899999999
$ws.Range("L132").Value = 10205.5716
$ws.Range("M132").Value = -6851.249899999999
$ws.Range("N132").Value = -15265.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 9624.875
$ws.Range("I13").Value = 2249.5
$ws.Range("J13").Value = 12083.333
$ws.Range("K13").Value = 2249.5
$ws.Range("L13").Value = 12083.333
$ws.Range("M13").Value = -2109.5
$ws.Range("N13").Value = -12363.333
$ws.Range("H74").Value = 82857.14
$ws.Range("J74").Value = 82857.14
$ws.Range("L74").Value = 82857.14
$ws.Range("N74").Value = -84853.14
$ws.Range("H77").Value = 82857.14
$ws.Range("J77").Value = 82857.14
$ws.Range("L77").Value = 248571.42
$ws.Range("N77").Value = -258555.42
$ws.Range("H116").Value = 243328
$ws.Range("J116").Value = 243328
$ws.Range("L116").Value = 243328
$ws.Range("N116").Value = -252506
$ws.Range("H131").Value = 69828.5
$ws.Range("J131").Value = 69828.5
$ws.Range("L131").Value = 69828.5
$ws.Range("N131").Value = -79908.5
$ws.Range("H136").Value = 2518.5156
$ws.Range("I136").Value = 1415.5476
$ws.Range("K136").Value = 4246.642800000001
$ws.Range("M136").Value = -1696.642800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 27780728
$ws.Range("I81").Value = 1400
$ws.Range("J81").Value = 55560056
$ws.Range("K81").Value = 2800
$ws.Range("L81").Value = 111120112
$ws.Range("M81").Value = -1739
$ws.Range("N81").Value = -111122234
$ws.Range("H84").Value = 27780728
$ws.Range("I84").Value = 1400
$ws.Range("J84").Value = 55560056
$ws.Range("K84").Value = 14000
$ws.Range("L84").Value = 555600560
$ws.Range("M84").Value = -8696
$ws.Range("N84").Value = -555611168
$ws.Range("H107").Value = 571.86487
$ws.Range("I107").Value = 580.80646
$ws.Range("K107").Value = 1742.41938
$ws.Range("M107").Value = 177.58062
$ws.Range("H136").Value = 2550.3125
$ws.Range("I136").Value = 2684.276
$ws.Range("K136").Value = 8052.828
$ws.Range("M136").Value = -5502.828

